$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated CasesTab Cypher query: split the WHERE/MATCH onto separate
# lines, compute Age/Weight with an integer-collapsing CASE expression,
# and add an ORDER BY / LIMIT 100 to the result set.
$query = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Mixed Breed']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$ws.Range("B2").Value = $query

# Column B now needs its own (wider) width since the query text grew;
# column C keeps the original shared width.
$ws.Columns.Item(2).ColumnWidth = 85.1667

# Row 2 grows taller to fit the longer wrapped query text.
$ws.Rows.Item(2).RowHeight = 345

# Scroll the view back to the top-left (A1) instead of A2.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
